$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'41.193.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -1.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.174.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -2.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'237.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -2.47%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  -2.11%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'70.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -5.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  +0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  -6.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'40.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Formula = "'  -3.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Formula = "'  -1.71%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Formula = "'  -5.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'2.496.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -2.17%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'13.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -2.27%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Formula = "'  -4.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'2.156.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  -2.70%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'40.970.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  -2.15%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'0.0000101"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -7.74%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'70.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  -2.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'5.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -4.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'10.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -4.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'225.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -1.71%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Formula = "'  -7.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Formula = "'  -0.03%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'10.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -5.75%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'3.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  -2.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Formula = "'  -3.44%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Formula = "'  +1.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'166.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +0.22%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'19.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  -3.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'30.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +4.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'0.0767"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -4.31%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  -9.38%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Formula = "'  -3.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Formula = "'  -9.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'4.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  -3.83%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'0.0286"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  -5.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  -4.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E41").Formula = "'  -4.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'60.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -7.67%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Formula = "'  -4.66%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'8.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  -4.59%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.0972"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -3.77%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'98.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Formula = "'  -2.53%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Formula = "'  -2.98%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'  -7.80%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  -2.99%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'2.374.71"
$ws.Range("D51").Style = "Normal"
